# Adds 20 new "DB" sheet rows for the 2023-05-30 (230530) recording
# session of subject BL-003 / trail 024, and marks the "DB" tab as the
# active tab (matches the author's "reading Time Counter from MC and
# writing in the .csv file" commit: new recordings appended to the DB).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DB")

# Recordings ("D" column) timestamps for 230530, in sheet order.
$timestamps = @(
    "230530_113523",
    "230530_113845",
    "230530_115123",
    "230530_120934",
    "230530_122720",
    "230530_124547",
    "230530_125056",
    "230530_125439",
    "230530_130747",
    "230530_131220",
    "230530_132329",
    "230530_132527",
    "230530_132832",
    "230530_133313",
    "230530_133611",
    "230530_134658",
    "230530_135617",
    "230530_135815",
    "230530_135940",
    "230530_140005"
)

# Matching "protocol" ("G" column) value for each row above.
$protocols = @(
    "Exp#5_14s",
    "Exp#5_14s",
    "Exp#1_60s",
    "Exp#5_14s",
    "Exp#1_60s",
    "Exp#1_60s",
    "Exp#1_60s",
    "Exp#1_60s",
    "Exp#1_60s",
    "Exp#1_60s",
    "Exp#1_60s",
    "Exp#1_60s",
    "Exp#1_60s",
    "Exp#1_60s",
    "Exp#1_60s",
    "Exp#1_60s",
    "Exp#1_60s",
    "Exp#5_14s",
    "Exp#5_14s",
    "Exp#5_14s"
)

$firstRow = 923
$lastRow = $firstRow + $timestamps.Length - 1

# Columns B..H are plain (non-numeric-looking) text, so a direct Value2
# assignment keeps them as shared-string text without Excel re-typing
# them as numbers.
$r = $firstRow
for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $ws.Range("B$r").Value2 = "BL-003"
    $ws.Range("C$r").Value2 = "BL-003_024_230530"
    $ws.Range("D$r").Value2 = $timestamps[$i]
    $ws.Range("E$r").Value2 = "1.1, 1.2"
    $ws.Range("F$r").Value2 = "above knee"
    $ws.Range("G$r").Value2 = $protocols[$i]
    $ws.Range("H$r").Value2 = "left leg"
    $r = $r + 1
}

# Column A holds the trail id "024" -- a numeric-looking string that a
# direct Value2 assignment would silently coerce to the number 24.
# Copy it (as a value-only paste) from the existing identical column-A
# block instead, so it lands as text without picking up a new number
# format / style.
$ws.Range("A883:A902").Copy() | Out-Null
$ws.Range("A$firstRow`:A$lastRow").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# Make "DB" the active tab (activeTab="3" in the saved workbook).
$ws.Activate()
